$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "Проход по дереву до null и создание arraylist<cells>"
$ws.Range("B11").Value = "Рекурсия : на вход массив cells  с текущей cell, наполнение одним из вариантов и вызов рекурсивный ."

$ws.Columns.Item(2).ColumnWidth = 49.166666666666664
$ws.Columns.Item(4).ColumnWidth = 25.666666666666668

$ws.Range("B12").Select()
